$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 171.5
$ws.Range("I8").Value = 171.5
$ws.Range("K8").Value = 514.5
$ws.Range("M8").Value = -375.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 289.3
$ws.Range("I53").Value = 79.3
$ws.Range("J53").Value = 499.3
$ws.Range("K53").Value = 79.3
$ws.Range("L53").Value = 499.3
$ws.Range("M53").Value = 557.7
$ws.Range("N53").Value = -1773.3

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 416.47058
$ws.Range("I55").Value = 489.85715
$ws.Range("J55").Value = 74
$ws.Range("K55").Value = 489.85715
$ws.Range("L55").Value = 74
$ws.Range("M55").Value = -275.85715
$ws.Range("N55").Value = -502

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 603.4666999999999
$ws.Range("I80").Value = 274.85715
$ws.Range("J80").Value = 891
$ws.Range("K80").Value = 824.5714499999999
$ws.Range("L80").Value = 2673
$ws.Range("M80").Value = 173.4285500000001
$ws.Range("N80").Value = -4669

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value = 603.4666999999999
$ws.Range("I83").Value = 274.85715
$ws.Range("J83").Value = 891
$ws.Range("K83").Value = 2473.71435
$ws.Range("L83").Value = 8019
$ws.Range("M83").Value = 2518.28565
$ws.Range("N83").Value = -18003

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 78549.266
$ws.Range("I135").Value = 49361
$ws.Range("J135").Value = 201140
$ws.Range("K135").Value = 444249
$ws.Range("L135").Value = 1810260
$ws.Range("M135").Value = -441714
$ws.Range("N135").Value = -1815330

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2087614.6
$ws.Range("I137").Value = 3705442.2
$ws.Range("J137").Value = 7550.524
$ws.Range("K137").Value = 11116326.6
$ws.Range("L137").Value = 22651.572
$ws.Range("M137").Value = -11113776.6
$ws.Range("N137").Value = -27751.572

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1785.4103
$ws.Range("I2").Value = 1692.1538
$ws.Range("J2").Value = 1971.9231
$ws.Range("K2").Value = 1692.1538
$ws.Range("L2").Value = 1971.9231
$ws.Range("M2").Value = -1579.1538
$ws.Range("N2").Value = -2197.9231

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1248.6364
$ws.Range("I45").Value = 579.4
$ws.Range("J45").Value = 1806.3334
$ws.Range("K45").Value = 579.4
$ws.Range("L45").Value = 1806.3334
$ws.Range("M45").Value = -202.4
$ws.Range("N45").Value = -2560.3334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 7291497
$ws.Range("I74").Value = 8155801.5
$ws.Range("J74").Value = 204202.8
$ws.Range("K74").Value = 8155801.5
$ws.Range("L74").Value = 204202.8
$ws.Range("M74").Value = -8154927.5
$ws.Range("N74").Value = -205950.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 7291497
$ws.Range("I77").Value = 8155801.5
$ws.Range("J77").Value = 204202.8
$ws.Range("K77").Value = 40779007.5
$ws.Range("L77").Value = 1021014
$ws.Range("M77").Value = -40774639.5
$ws.Range("N77").Value = -1029750

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1785.4103
$ws.Range("I116").Value = 1692.1538
$ws.Range("J116").Value = 1971.9231
$ws.Range("K116").Value = 1692.1538
$ws.Range("L116").Value = 1971.9231
$ws.Range("M116").Value = 601.8462
$ws.Range("N116").Value = -6559.9231

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1785.4103
$ws.Range("I3").Value = 1692.1538
$ws.Range("J3").Value = 1971.9231
$ws.Range("K3").Value = 1692.1538
$ws.Range("L3").Value = 1971.9231
$ws.Range("M3").Value = -1578.1538
$ws.Range("N3").Value = -2199.9231

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 439.83334
$ws.Range("I22").Value = 439.83334
$ws.Range("K22").Value = 439.83334
$ws.Range("M22").Value = -266.83334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 15739.889
$ws.Range("I82").Value = 6113.1665
$ws.Range("J82").Value = 34993.332
$ws.Range("K82").Value = 6113.1665
$ws.Range("L82").Value = 34993.332
$ws.Range("M82").Value = -5730.1665
$ws.Range("N82").Value = -35759.332

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H85").Value = 15739.889
$ws.Range("I85").Value = 6113.1665
$ws.Range("J85").Value = 34993.332
$ws.Range("K85").Value = 6113.1665
$ws.Range("L85").Value = 34993.332
$ws.Range("M85").Value = -4787.1665
$ws.Range("N85").Value = -37645.332

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 34016.414
$ws.Range("I31").Value = 23576.32
$ws.Range("J31").Value = 43337.93
$ws.Range("K31").Value = 23576.32
$ws.Range("L31").Value = 43337.93
$ws.Range("M31").Value = -23281.32
$ws.Range("N31").Value = -43927.93

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 34016.414
$ws.Range("I34").Value = 23576.32
$ws.Range("J34").Value = 43337.93
$ws.Range("K34").Value = 23576.32
$ws.Range("L34").Value = 43337.93
$ws.Range("M34").Value = -23374.32
$ws.Range("N34").Value = -43741.93

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 22223652
$ws.Range("I58").Value = 31251352
$ws.Range("J58").Value = 1623.1538
$ws.Range("K58").Value = 31251352
$ws.Range("L58").Value = 1623.1538
$ws.Range("M58").Value = -31251149
$ws.Range("N58").Value = -2029.1538

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1520.1875
$ws.Range("I122").Value = 1300.4814
$ws.Range("J122").Value = 2706.6
$ws.Range("K122").Value = 3901.4442
$ws.Range("L122").Value = 8119.799999999999
$ws.Range("M122").Value = -1451.4442
$ws.Range("N122").Value = -13019.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 22223652
$ws.Range("I136").Value = 31251352
$ws.Range("J136").Value = 1623.1538
$ws.Range("K136").Value = 93754056
$ws.Range("L136").Value = 4869.4614
$ws.Range("M136").Value = -93751506
$ws.Range("N136").Value = -9969.4614

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 928.1142599999999
$ws.Range("I131").Value = 351.8
$ws.Range("J131").Value = 1024.1666
$ws.Range("K131").Value = 1055.4
$ws.Range("L131").Value = 3072.4998
$ws.Range("M131").Value = 3984.6
$ws.Range("N131").Value = -13152.4998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 2315.6216
$ws.Range("I132").Value = 1313.0416
$ws.Range("K132").Value = 11817.3744
$ws.Range("M132").Value = -9287.374400000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3328.963
$ws.Range("I80").Value = 2673.8
$ws.Range("J80").Value = 4147.9165
$ws.Range("K80").Value = 2673.8
$ws.Range("L80").Value = 4147.9165
$ws.Range("M80").Value = -1675.8
$ws.Range("N80").Value = -6143.9165

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 3328.963
$ws.Range("I83").Value = 2673.8
$ws.Range("J83").Value = 4147.9165
$ws.Range("K83").Value = 13369
$ws.Range("L83").Value = 20739.5825
$ws.Range("M83").Value = -8377
$ws.Range("N83").Value = -30723.5825

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1784.8
$ws.Range("I102").Value = 1841.3334
$ws.Range("J102").Value = 1700
$ws.Range("K102").Value = 1841.3334
$ws.Range("L102").Value = 1700
$ws.Range("M102").Value = -219.3334
$ws.Range("N102").Value = -4944

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 130.46666
$ws.Range("I55").Value = 136.21428
$ws.Range("J55").Value = 50
$ws.Range("K55").Value = 136.21428
$ws.Range("L55").Value = 50
$ws.Range("M55").Value = 36.78572
$ws.Range("N55").Value = -396

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3264.347
$ws.Range("I122").Value = 2595.1875
$ws.Range("J122").Value = 3588.7878
$ws.Range("K122").Value = 7785.5625
$ws.Range("L122").Value = 10766.3634
$ws.Range("M122").Value = -5335.5625
$ws.Range("N122").Value = -15666.3634

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 7570.857
$ws.Range("I41").Value = 8080.6665
$ws.Range("J41").Value = 7188.5
$ws.Range("K41").Value = 8080.6665
$ws.Range("L41").Value = 7188.5
$ws.Range("M41").Value = -7690.6665
$ws.Range("N41").Value = -7968.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 914.8
$ws.Range("I126").Value = 755.1539
$ws.Range("K126").Value = 2265.4617
$ws.Range("M126").Value = 204.5383000000002
